$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C values (changed order of s)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 6

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("C6").Select()

# Reflect the updated window geometry from the workbook view (best effort)
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 5960
    $win.Top = 1460
    $win.Width = 10000
    $win.Height = 15880
} catch {}
